$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "01tq0000001jgnm"
$ws.Range("C10").Select()
